$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.787.07'
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.813.03'
$ws.Range("E3").Value = '  +1.58%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '351.03'
$ws.Range("E5").Value = '  -0.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.18'
$ws.Range("E6").Value = '  +5.37%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.561'
$ws.Range("E7").Value = '  +2.17%  '

$ws.Range("E9").Value = '  +6.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.22'
$ws.Range("E10").Value = '  +2.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0850'
$ws.Range("E12").Value = '  +2.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.12'
$ws.Range("E13").Value = '  +0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.79'
$ws.Range("E14").Value = '  +3.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.253.63'
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.974'
$ws.Range("E16").Value = '  +5.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.821.76'
$ws.Range("E17").Value = '  +2.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.841.26'
$ws.Range("E18").Value = '  +1.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.38'
$ws.Range("E19").Value = '  +9.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.61'
$ws.Range("E20").Value = '  -0.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.48'
$ws.Range("E21").Value = '  +3.15%  '

$ws.Range("E22").Value = '  +1.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.63'
$ws.Range("E23").Value = '  +1.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.46'
$ws.Range("E24").Value = '  +1.09%  '

$ws.Range("E25").Value = '  +1.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.22'
$ws.Range("E26").Value = '  +1.11%  '

$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("E28").Value = '  +0.10%  '

$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.78'
$ws.Range("E29").Value = '  +10.64%  '

$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.53'
$ws.Range("E30").Value = '  +3.66%  '

$ws.Range("E31").Value = '  -1.36%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.23'
$ws.Range("E32").Value = '  +2.81%  '

$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '52.72'
$ws.Range("E33").Value = '  +1.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.76'
$ws.Range("E34").Value = '  +4.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0904'
$ws.Range("E35").Value = '  +9.67%  '

$ws.Range("E36").Value = '  +2.48%  '

$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.00'
$ws.Range("E38").Value = '  +4.76%  '

$ws.Range("E39").Value = '  +2.19%  '

$ws.Range("E40").Value = '  +2.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.57'
$ws.Range("E41").Value = '  +2.65%  '

$ws.Range("E42").Value = '  +1.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '121.60'

$ws.Range("E44").Value = '  +1.89%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.08'
$ws.Range("E45").Value = '  +0.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.52'
$ws.Range("E46").Value = '  +8.87%  '

$ws.Range("E47").Value = '  +9.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.146.90'
$ws.Range("E48").Value = '  +2.55%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.989'
$ws.Range("E49").Value = '  +8.64%  '

$ws.Range("E50").Value = '  +19.17%  '

$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0321'
$ws.Range("E51").Value = '  +15.37%  '
